$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 3 new rows in the "Casing" section (rows 26-28), pushing the
#    "Material"/"Tools" sections down by 3 rows.
$ws.Rows("26:28").Insert()

# 2) Add the new trailing row under "Tools" section (now row 41): Atmel ICE
#    programmer. Doing this before filling in rows 26-28 makes it claim
#    shared-string index 38, matching the author's original edit order
#    captured in the diff.
$ws.Range("A41").Value = "Atmel ICE programmer"
$ws.Range("C41").Value = "Need"

# 3) Fill in the three newly inserted rows.
$ws.Range("A26").Value = "Hinges"
$ws.Range("B26").Value = 2
$ws.Range("C26").Value = "Need"

$ws.Range("A27").Value = "Latch"
$ws.Range("B27").Value = 1
$ws.Range("C27").Value = "Need"

$ws.Range("A28").Value = "Rubber feet"
$ws.Range("B28").Value = 4
$ws.Range("C28").Value = "Need"

# 4) Column A needs to widen (auto-fit) to accommodate the longest new entry
#    ("Atmel ICE programmer"); only touch column A so B/C keep their bestFit.
$ws.Columns("A").ColumnWidth = 18.8333333

# 5) Move the active selection to D28 (the frozen pane's scroll position
#    naturally resets to the top after the row insert above).
$ws.Range("D28").Select()
